# Update attendance ("想去人数") and minimum ticket price ("最低票价")
# figures for both the "展览" (sheet1) and "全部类型" (sheet4) worksheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 1102
    $ws.Range("G2").Value = 68

    $ws.Range("F5").Value = 8772

    $ws.Range("F8").Value = 649

    $ws.Range("F10").Value = 158

    $ws.Range("F13").Value = 3622

    $ws.Range("F15").Value = 367

    $ws.Range("F16").Value = 82

    $ws.Range("F17").Value = 1906

    $ws.Range("F18").Value = 150

    $ws.Range("F20").Value = 313

    $ws.Range("F21").Value = 208

    $ws.Range("F22").Value = 2408
}
